# Daily attendance processing - normalize "Recorded By" (column G) values
# so that the literal token "System" (capitalized, case-sensitive) always
# appears first in the comma-separated list of recorders, preserving the
# relative order of all other tokens (including a lowercase "system"
# entry, if present, which is left where it is).

function Test-IsExactSystem($s) {
    if ($s -ne "System") { return $false }
    $code0 = [int][char]$s[0]
    return ($code0 -eq 83)   # 'S' (uppercase) ASCII code
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $systemIndex = -1
        for ($i = 0; $i -lt $trimmed.Count; $i++) {
            if (Test-IsExactSystem $trimmed[$i]) {
                $systemIndex = $i
                break
            }
        }

        if ($systemIndex -gt 0) {
            $reordered = @("System")
            for ($i = 0; $i -lt $trimmed.Count; $i++) {
                if ($i -ne $systemIndex) {
                    $reordered += $trimmed[$i]
                }
            }
            $cell.Value2 = ($reordered -join ", ")
        }
    }
}
